# Naive component forecaster bug fix - Presentation state 11.02.
# A new "Q0" forecast-error observation is prepended to each row's series
# (columns B..K hold the 10 most-recent quarters, B = most recent).
# Existing values shift one column to the right; once a row's window is
# full (columns B..K all populated) the oldest observation (column K)
# falls off the back of the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newvals = @{
  2  = 0.9564081874156993
  3  = -4.157449276732949
  4  = 1.546611864454844
  5  = 1.156631887942306
  6  = -1.025188112727922
  7  = 0.08364543516793629
  8  = -0.1538585523806955
  9  = 0.7495351060200912
  10 = 0.03849281619118239
  11 = -0.2590580299438133
  12 = 0.01855976243503714
  13 = 0.1467044301255134
  14 = -0.1819613811903656
  15 = 0.4718454808444464
  16 = -0.08594117411414147
  17 = -0.07695400962807622
  18 = -0.5068991247689255
  19 = 0.6215838649243215
  20 = -0.2766911554241067
}

for ($row = 2; $row -le 20; $row++) {
    $val = $newvals[$row]

    # Find the last populated column in B..K (2..11) for this row.
    $lastCol = 1
    for ($c = 2; $c -le 11; $c++) {
        if ($ws.Cells.Item($row, $c).Value() -ne $null) {
            $lastCol = $c
        }
    }

    # Shift existing values one column to the right (from the end
    # backwards so we don't overwrite values before they're copied).
    # Anything that would land past column K (11) drops off.
    for ($c = $lastCol; $c -ge 2; $c--) {
        $srcVal = $ws.Cells.Item($row, $c).Value()
        $destCol = $c + 1
        if ($destCol -le 11) {
            $ws.Cells.Item($row, $destCol).Value = $srcVal
        }
    }

    # Prepend the new observation in column B.
    $ws.Cells.Item($row, 2).Value = $val
}
